# Update column F ("dSF") values for a set of rows in Sheet1.
# These values were "repulled" from source data; only specific rows changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -8
    3  = 0
    4  = -6
    5  = 3
    8  = -3
    9  = -4
    10 = -1
    11 = -4
    12 = 0
    16 = -5
    22 = 0
    24 = 2
    27 = 0
    31 = -3
    32 = 5
    35 = -1
    36 = 1
    38 = -7
    40 = -1
    44 = -3
    46 = -1
    47 = -4
    48 = 3
    56 = 5
    57 = -5
    58 = -2
    62 = -2
    64 = 1
    65 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
